# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (want-to-go count) figures in column F across the
# four worksheets: 展览 (Exhibitions), 演出 (Performances), 本地生活 (Local
# Life) and 全部类型 (All Types, the combined/merged view of the other three).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2621    # was 2612
$ws.Range("F5").Value = 918     # was 917
$ws.Range("F7").Value = 1925    # was 1923
$ws.Range("F8").Value = 1770    # was 1763
$ws.Range("F9").Value = 201     # was 200
$ws.Range("F11").Value = 2428   # was 2423
$ws.Range("F12").Value = 524    # was 522
$ws.Range("F13").Value = 203    # was 197
$ws.Range("F18").Value = 8919   # was 8894
$ws.Range("F20").Value = 6937   # was 6920
$ws.Range("F21").Value = 11320  # was 11279
$ws.Range("F24").Value = 223    # was 222
$ws.Range("F25").Value = 309    # was 311
$ws.Range("F27").Value = 2475   # was 2468
$ws.Range("F28").Value = 218    # was 213
$ws.Range("F30").Value = 2384   # was 2362
$ws.Range("F31").Value = 578    # was 555
$ws.Range("F32").Value = 36     # was 34
$ws.Range("F33").Value = 4479   # was 4477
$ws.Range("F34").Value = 742    # was 709
$ws.Range("F35").Value = 321    # was 315
$ws.Range("F36").Value = 28     # was 26
$ws.Range("F37").Value = 486    # was 481

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F23").Value = 6      # was 5

# 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 127     # was 121

# 全部类型 (All Types - combined view)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 2621    # was 2612
$ws.Range("F8").Value = 918     # was 917
$ws.Range("F10").Value = 1925   # was 1923
$ws.Range("F12").Value = 1770   # was 1763
$ws.Range("F14").Value = 201    # was 200
$ws.Range("F15").Value = 2428   # was 2423
$ws.Range("F17").Value = 524    # was 522
$ws.Range("F18").Value = 203    # was 197
$ws.Range("F23").Value = 8920   # was 8894
$ws.Range("F25").Value = 6937   # was 6920
$ws.Range("F26").Value = 11320  # was 11279
$ws.Range("F29").Value = 223    # was 222
$ws.Range("F30").Value = 309    # was 310
$ws.Range("F36").Value = 220    # was 213
$ws.Range("F38").Value = 36     # was 34
$ws.Range("F39").Value = 4479   # was 4477
$ws.Range("F46").Value = 486    # was 481
$ws.Range("F49").Value = 6      # was 5
